$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("United States of America, Government of", 14049432),
    @("United Arab Emirates, Government of", 3003000),
    @("Germany, Government of", 710201),
    @("Canada, Government of", 735864),
    @("Sweden, Government of", 921037),
    @("Central Emergency Response Fund", 5167962),
    @("US Fund for UNICEF", 139607),
    @("European Commission's Humanitarian Aid and Civil Protection Department", 4047077),
    @("Education Cannot Wait Fund", 250000),
    @("Italy, Government of", 21277),
    @("United Nations Population Fund", 300000),
    @("Spain, Government of", 313589),
    @("Switzerland, Government of", 730124),
    @("Argentina, Government of", 0),
    @("Paraguay, Government of", 0),
    @("Venezuela, Government of", 0),
    @("Bolivia, Government of", 0),
    @("Brazil, Government of", 0),
    @("Chile, Government of", 0),
    @("Colombia, Government of", 0),
    @("Ecuador, Government of", 0),
    @("European Commission Civil Protection Mechanism", 0),
    @("Korea, Republic of, Government of", 300000),
    @("Qatar, Government of", 1000000)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
